$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E and G values (col F "Win" is unchanged) - regenerated sval data
# that filters save games, per commit message.
$data = @{
    2  = @{ B = 1.505614041169197;   C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 4.371470058157054 }
    3  = @{ B = 0.06328177979961902; C = 0.3375848360084654;  D = 0.7127328510149897;  E = 0.4998867070740569; G = 1.613486173897131 }
    4  = @{ B = 0.1554434735375247;  C = 0.3375848360084654;  D = 3.082599426703578;   E = 0.4998867070740569; G = 4.075514443323626 }
    5  = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812;  E = 0.4998867070740569; G = 5.488907176552729 }
    6  = @{ B = 0.3464964993005633;  C = 0.3375848360084654;  D = 0.1529057820181812;  E = 0.4998867070740569; G = 1.336873824401267 }
    7  = @{ B = 1.505614041169197;   C = 1.65323645889881;    D = 3.082599426703578;   E = 0.4998867070740569; G = 6.741336633845642 }
    8  = @{ B = 0.06328177979961902; C = 0.05231270169004087; D = 0.7127328510149897;  E = 0.4998867070740569; G = 1.328214039578707 }
    9  = @{ B = 0.7287194209349384;  C = 0.004309184025731883;D = 0.1529057820181812;  E = 0.4998867070740569; G = 1.385821094052908 }
    10 = @{ B = 0.02258322285507441; C = 0.3375848360084654;  D = 0.7127328510149897;  E = 0.4998867070740569; G = 1.572787616952587 }
    11 = @{ B = 0.3464964993005633;  C = 1.65323645889881;    D = 0.1529057820181812;  E = 6.48142807727062;    G = 8.634066817488176 }
    12 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 3.082599426703578;   E = 0.4998867070740569; G = 8.418600821238126 }
    13 = @{ B = 1.505614041169197;   C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 4.371470058157054 }
    14 = @{ B = 0.06328177979961902; C = 0.05231270169004087; D = 3.082599426703578;   E = 0.4998867070740569; G = 3.698080615267295 }
    15 = @{ B = 0.7287194209349384;  C = 1.65323645889881;    D = 0.7127328510149897;  E = 0.4998867070740569; G = 3.594575437922795 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
